$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'313.94"
$ws.Range("E2").Value = "'2.01%"
$ws.Range("D3").Value = "'40.90"
$ws.Range("E3").Value = "'-0.34%"
$ws.Range("E4").Value = "'-1.54%"
$ws.Range("D5").Value = "'0.07594"
$ws.Range("E5").Value = "'-0.97%"
$ws.Range("D6").Value = "'4.329"
$ws.Range("E6").Value = "'0.41%"
$ws.Range("D7").Value = "'1.680"
$ws.Range("E7").Value = "'2.19%"
$ws.Range("D8").Value = "'0.9301"
$ws.Range("E8").Value = "'1.58%"
$ws.Range("D10").Value = "'0.1198"
$ws.Range("E10").Value = "'-3.90%"
$ws.Range("D11").Value = "'0.1815"
$ws.Range("E11").Value = "'-0.71%"
$ws.Range("D12").Value = "'0.09052"
$ws.Range("E12").Value = "'-0.84%"
$ws.Range("D13").Value = "'0.04139"
$ws.Range("E13").Value = "'-2.56%"
$ws.Range("E14").Value = "'0.39%"
$ws.Range("D15").Value = "'0.001292"
$ws.Range("E15").Value = "'2.45%"
$ws.Range("D16").Value = "'0.005833"
$ws.Range("E16").Value = "'1.21%"
$ws.Range("D18").Value = "'3.330"
$ws.Range("E18").Value = "'-0.52%"
$ws.Range("D19").Value = "'0.3358"
$ws.Range("E19").Value = "'0.66%"
$ws.Range("D20").Value = "'7.614"
$ws.Range("E20").Value = "'4.10%"
$ws.Range("E21").Value = "'-3.08%"
$ws.Range("D22").Value = "'0.2838"
$ws.Range("E22").Value = "'-1.95%"
$ws.Range("D23").Value = "'0.04028"
$ws.Range("E23").Value = "'-1.13%"
$ws.Range("E24").Value = "'1.13%"
$ws.Range("D25").Value = "'0.003971"
$ws.Range("E25").Value = "'-7.24%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'2.26%"
$ws.Range("D38").Value = "'0.02412"
$ws.Range("E38").Value = "'-2.45%"
$ws.Range("D39").Value = "'0.05167"
$ws.Range("E39").Value = "'-2.32%"
$ws.Range("D40").Value = "'0.007746"
$ws.Range("E40").Value = "'-1.23%"
$ws.Range("D41").Value = "'0.1301"
$ws.Range("E41").Value = "'-0.96%"
$ws.Range("D42").Value = "'0.007609"
$ws.Range("E42").Value = "'10.67%"
$ws.Range("E43").Value = "'72.58%"
$ws.Range("D44").Value = "'0.008583"
$ws.Range("E44").Value = "'12.47%"
$ws.Range("D45").Value = "'0.3386"
$ws.Range("E45").Value = "'10.67%"
$ws.Range("D46").Value = "'0.00006588"
$ws.Range("E46").Value = "'-2.02%"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("D48").Value = "'0.2686"
$ws.Range("E48").Value = "'58.14%"
$ws.Range("D49").Value = "'0.004203"
$ws.Range("E49").Value = "'35.35%"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("E51").Value = "'-0.11%"
